# ITSADSSD-21887 - Placeholder for git and code changes for send email
#
# Settings sheet: add PBI input-file / email-exchange / subject / archive
# folder settings, renumbering several existing rows; Assets sheet: add
# UQ SMTP server/port asset rows; switch the active tab to "Assets".

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("Settings")
$assets   = $wb.Worksheets.Item("Assets")

# ---------------------------------------------------------------------
# Settings sheet: clear the old row6 (EmailExchangeServer moves to row7,
# leaving row6 blank) before rewriting everything below row 2.
# ---------------------------------------------------------------------
$settings.Range("A6:C6").Clear() | Out-Null

# Row 3: Folder_Location
$settings.Range("A3").Value = "Folder_Location"
$settings.Range("B3").Value = "Data"
$settings.Range("C3").Value = "Folder to save input file "

# Row 4: File_Name
$settings.Range("A4").Value = "File_Name"
$settings.Range("B4").Value = "EXTRACT_FOR_ROBOT.csv"
$settings.Range("C4").Value = "File name for input, received from PBI team"

# Row 5: logF_BusinessProcessName (was row 4; value changes Framework -> PBI_LocationCheck)
$settings.Range("A5").Value = "logF_BusinessProcessName"
$settings.Range("B5").Value = "PBI_LocationCheck"
$settings.Range("C5").Value = "This is a logging field which allows you to group the log data of two or more subprocesses under the same business process name"

# Row 7: EmailExchangeServer (was row 6), keeps the Hyperlink style
$settings.Range("A7").Value = "EmailExchangeServer"
$settings.Range("B7").Value = "https://outlook.office365.com/EWS/Exchange.asmx"
$settings.Range("C7").Value = "Exchange server web address to retrive emails "

# Row 8: NumberOfEmails (was row 7), keeps the left-aligned numeric style
$settings.Range("A8").Value = "NumberOfEmails"
$settings.Range("B8").Value = 10
$settings.Range("B8").HorizontalAlignment = -4131
$settings.Range("C8").Value = "Number of emails to check (Top attribute in GetExchangeEmailessage activity) "

# Row 9: Subject_Prefix (was row 8)
$settings.Range("A9").Value = "Subject_Prefix"
$settings.Range("B9").Value = "Location VEVO check"
$settings.Range("C9").Value = "Subject for email received from PBI team ro perform check"

# Row 10: FolderName (new)
$settings.Range("A10").Value = "FolderName"
$settings.Range("B10").Value = "PBI_Archive"
$settings.Range("C10").Value = "Folder to move email after processing "

# Row 12: Immi_Website (was row 10), keeps the Hyperlink style
$settings.Range("A12").Value = "Immi_Website"
$settings.Range("B12").Value = "https://online.immi.gov.au/lusc/login"

# Rebuild the hyperlinks collection against the new anchor cells, then
# reapply the named "Hyperlink" style so the cells land back on the same
# shared style slot the workbook already uses (Hyperlinks.Add otherwise
# mints its own near-duplicate style).
$settings.Hyperlinks.Delete() | Out-Null
$settings.Hyperlinks.Add($settings.Range("B7"), "https://outlook.office365.com/EWS/Exchange.asmx") | Out-Null
$settings.Hyperlinks.Add($settings.Range("B12"), "https://online.immi.gov.au/lusc/login") | Out-Null
$settings.Range("B7").Style = "Hyperlink"
$settings.Range("B12").Style = "Hyperlink"

# Extend the sheet's used range down to row 999 (two extra trailing blank
# rows), matching the refreshed dimension in the target workbook.
$settings.Range("A998").Font.Bold = $false
$settings.Range("A999").Font.Bold = $false
$settings.Range("A998:A999").ClearContents() | Out-Null

# ---------------------------------------------------------------------
# Assets sheet: add descriptions for the existing credential/asset rows
# and append the new UQ SMTP server/port assets.
# ---------------------------------------------------------------------
$assets.Range("C2").Value = "Credentials for SSO account  RPA00001"
$assets.Range("C3").Value = "Credentials for Immi website login"
$assets.Range("C4").Value = "PBI team email address "

$assets.Range("A5").Value = "UQ_SMTP_SERVER"
$assets.Range("B5").Value = "UQ_SMTP_SERVER"
$assets.Range("C5").Value = "UQ SMTP server address"

$assets.Range("A6").Value = "UQ_SMTP_PORT"
$assets.Range("B6").Value = "UQ_SMTP_PORT"
$assets.Range("C6").Value = "UQ SMTP port number "

# ---------------------------------------------------------------------
# Selection / active tab: Assets becomes the active sheet (C18 selected);
# Settings keeps A3 selected but is no longer the active tab.
# ---------------------------------------------------------------------
$settings.Range("A3").Select() | Out-Null
$assets.Activate() | Out-Null
$assets.Range("C18").Select() | Out-Null
